$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = -3.377448218355161
$ws.Cells.Item(2, 3).Value = 2.15496247269007
$ws.Cells.Item(2, 4).Value = 7.434495881238234
$ws.Cells.Item(3, 2).Value = 0.2448820992536849
$ws.Cells.Item(3, 3).Value = -0.1721910504140167
$ws.Cells.Item(3, 4).Value = -1.008236814482644
$ws.Cells.Item(4, 2).Value = 5.624925926329238
$ws.Cells.Item(4, 3).Value = 0.6459332894185987
$ws.Cells.Item(4, 4).Value = 6.267294698459325
$ws.Cells.Item(5, 2).Value = 6.16826137189288
$ws.Cells.Item(5, 3).Value = -6.935153294182605
$ws.Cells.Item(5, 4).Value = 9.755477202390939
$ws.Cells.Item(6, 2).Value = -1.494625744833378
$ws.Cells.Item(6, 3).Value = -6.157851153201799
$ws.Cells.Item(6, 4).Value = 8.009006717074318
$ws.Cells.Item(7, 2).Value = -0.2554344193826941
$ws.Cells.Item(7, 3).Value = -4.932397146504464
$ws.Cells.Item(7, 4).Value = 2.809173804671983
$ws.Cells.Item(8, 2).Value = -0.8961894651313584
$ws.Cells.Item(8, 3).Value = -4.264270422584227
$ws.Cells.Item(8, 4).Value = 0.1319837197746532
$ws.Cells.Item(9, 2).Value = 4.683070112298715
$ws.Cells.Item(9, 3).Value = -1.149109367878443
$ws.Cells.Item(9, 4).Value = 11.26013612946959
$ws.Cells.Item(10, 2).Value = -10.69158489251918
$ws.Cells.Item(10, 3).Value = -5.600889391489416
$ws.Cells.Item(10, 4).Value = -5.900781225340602
$ws.Cells.Item(11, 2).Value = -6.373529693431146
$ws.Cells.Item(11, 3).Value = 9.024371965595002
$ws.Cells.Item(11, 4).Value = -7.350727472305019
$ws.Cells.Item(12, 2).Value = 0.04395139504043133
$ws.Cells.Item(12, 3).Value = 7.481601158193651
$ws.Cells.Item(12, 4).Value = -10.82795411725803
$ws.Cells.Item(13, 2).Value = -2.445450002465022
$ws.Cells.Item(13, 3).Value = 2.990192558263849
$ws.Cells.Item(13, 4).Value = -3.269703462328233
